$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values for rows 3, 5, 6, and 20 per the FlashScore data refresh.

# Row 3
$ws.Range("G3").Value = 1.95
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 4.1
$ws.Range("J3").Value = 2.75
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 8
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 17
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 67
$ws.Range("AG3").Value = 9
$ws.Range("AH3").Value = 19
$ws.Range("AJ3").Value = 41
$ws.Range("AU3").Value = 9
$ws.Range("AX3").Value = 23
$ws.Range("AY3").Value = 34
$ws.Range("AZ3").Value = 81
$ws.Range("BA3").Value = 126

# Row 5
$ws.Range("G5").Value = 1.83
$ws.Range("I5").Value = 4.1
$ws.Range("J5").Value = 2.6
$ws.Range("L5").Value = 5
$ws.Range("W5").Value = 5.5
$ws.Range("X5").Value = 7.5
$ws.Range("Y5").Value = 9
$ws.Range("Z5").Value = 15
$ws.Range("AH5").Value = 19
$ws.Range("AI5").Value = 15
$ws.Range("AJ5").Value = 51
$ws.Range("AK5").Value = 41
$ws.Range("AL5").Value = 51
$ws.Range("AN5").Value = 3.6
$ws.Range("AO5").Value = 10
$ws.Range("AU5").Value = 9.5
$ws.Range("AX5").Value = 26
$ws.Range("AY5").Value = 41
$ws.Range("AZ5").Value = 101
$ws.Range("BA5").Value = 151

# Row 6
$ws.Range("J6").Value = 2.9
$ws.Range("L6").Value = 3.3
$ws.Range("M6").Value = 10.9
$ws.Range("N6").Value = 1.04
$ws.Range("P6").Value = 3.25
$ws.Range("Q6").Value = 1.78
$ws.Range("R6").Value = 1.93
$ws.Range("U6").Value = 1.6
$ws.Range("V6").Value = 2.07
$ws.Range("W6").Value = 9
$ws.Range("X6").Value = 12.5
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 25
$ws.Range("AA6").Value = 18
$ws.Range("AB6").Value = 25
$ws.Range("AC6").Value = 11.25
$ws.Range("AE6").Value = 12.5
$ws.Range("AG6").Value = 9.75
$ws.Range("AO6").Value = 12
$ws.Range("AP6").Value = 19
$ws.Range("AR6").Value = 75
$ws.Range("AT6").Value = 2.6
$ws.Range("AU6").Value = 6.7
$ws.Range("AY6").Value = 21
$ws.Range("AZ6").Value = 65
$ws.Range("BA6").Value = 100

# Row 20
$ws.Range("G20").Value = 2.63
$ws.Range("I20").Value = 3.1
$ws.Range("J20").Value = 3.25
$ws.Range("L20").Value = 3.6
$ws.Range("M20").Value = 1.1
$ws.Range("N20").Value = 7
$ws.Range("Q20").Value = 2.3
$ws.Range("R20").Value = 1.6
$ws.Range("X20").Value = 12
$ws.Range("Y20").Value = 11
$ws.Range("Z20").Value = 26
$ws.Range("AA20").Value = 23
$ws.Range("AC20").Value = 7
$ws.Range("AG20").Value = 8.5
$ws.Range("AJ20").Value = 29
$ws.Range("AK20").Value = 26